$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '68.120.61'
Set-TextValue 'E2' '  +0.50%  '

Set-TextValue 'D3' '2.543.15'
Set-TextValue 'E3' '  +0.43%  '

Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.09%  '

Set-TextValue 'D5' '592.29'
Set-TextValue 'E5' '  +0.13%  '

Set-TextValue 'D6' '174.24'
Set-TextValue 'E6' '  -0.74%  '

Set-TextValue 'E7' '  -0.10%  '

Set-TextValue 'D8' '0.525'
Set-TextValue 'E8' '  -1.12%  '

Set-TextValue 'D9' '2.539.18'
Set-TextValue 'E9' '  +0.30%  '

Set-TextValue 'E10' '  -2.02%  '

Set-TextValue 'E11' '  +1.81%  '

Set-TextValue 'D12' '0.345'
Set-TextValue 'E12' '  +0.11%  '

Set-TextValue 'D13' '5.04'
Set-TextValue 'E13' '  -2.57%  '

Set-TextValue 'B14' 'Avalanche'
Set-TextValue 'C14' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D14' '26.47'
Set-TextValue 'E14' '  -1.05%  '

Set-TextValue 'B15' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C15' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D15' '3.017.98'
Set-TextValue 'E15' '  +0.71%  '

Set-TextValue 'B16' 'ShibaInu'
Set-TextValue 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.0000177'
Set-TextValue 'E16' '  -0.88%  '

Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '67.959.76'
Set-TextValue 'E17' '  +0.69%  '

Set-TextValue 'B18' 'Binance-PegBSC-USD'
Set-TextValue 'C18' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D18' '2.34'
Set-TextValue 'E18' '  +134.66%  '

Set-TextValue 'D19' '2.523.14'
Set-TextValue 'E19' '  +0.01%  '

Set-TextValue 'D20' '11.84'
Set-TextValue 'E20' '  +3.49%  '

Set-TextValue 'D21' '8.05'
Set-TextValue 'E21' '  -1.26%  '

Set-TextValue 'D22' '371.33'
Set-TextValue 'E22' '  +3.40%  '

Set-TextValue 'D23' '4.15'
Set-TextValue 'E23' '  -1.02%  '

Set-TextValue 'D24' '4.58'
Set-TextValue 'E24' '  -1.43%  '

Set-TextValue 'D25' '72.07'
Set-TextValue 'E25' '  +2.89%  '

Set-TextValue 'E26' '  -0.01%  '

Set-TextValue 'D27' '1.92'
Set-TextValue 'E27' '  -4.01%  '

Set-TextValue 'D28' '10.01'
Set-TextValue 'E28' '  -2.10%  '

Set-TextValue 'D30' '0.0₃0970'
Set-TextValue 'E30' '  -2.32%  '

Set-TextValue 'D31' '541.29'
Set-TextValue 'E31' '  -1.93%  '

Set-TextValue 'D32' '8.35'
Set-TextValue 'E32' '  +0.77%  '

Set-TextValue 'E33' '  -2.30%  '

Set-TextValue 'D34' '1.87'
Set-TextValue 'E34' '  +0.44%  '

Set-TextValue 'E35' '  -1.23%  '

Set-TextValue 'E36' '  -0.11%  '

Set-TextValue 'D37' '160.20'
Set-TextValue 'E37' '  +1.43%  '

Set-TextValue 'E38' '  -2.09%  '

Set-TextValue 'D39' '19.30'
Set-TextValue 'E39' '  +2.76%  '

Set-TextValue 'D40' '18.62'
Set-TextValue 'E40' '  +0.08%  '

Set-TextValue 'B41' 'RenderToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D41' '5.15'
Set-TextValue 'E41' '  -0.21%  '

Set-TextValue 'B42' 'PolygonEcosystemToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D42' '0.350'
Set-TextValue 'E42' '  -1.58%  '

Set-TextValue 'D43' '1.78'
Set-TextValue 'E43' '  -1.66%  '

Set-TextValue 'D44' '2.54'
Set-TextValue 'E44' '  -1.30%  '

Set-TextValue 'D45' '0.999'
Set-TextValue 'E45' '  -0.16%  '

Set-TextValue 'D46' '39.38'
Set-TextValue 'E46' '  -1.16%  '

Set-TextValue 'D47' '0.0₆0284'
Set-TextValue 'E47' '  +2.30%  '

Set-TextValue 'D48' '148.49'
Set-TextValue 'E48' '  -0.66%  '

Set-TextValue 'D49' '3.71'
Set-TextValue 'E49' '  +0.18%  '

Set-TextValue 'D50' '0.552'
Set-TextValue 'E50' '  -1.54%  '

Set-TextValue 'D51' '1.72'
Set-TextValue 'E51' '  +1.30%  '
